# Auto-generated edit script applying value updates to Sheets/Typhon_Profits.xlsx
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 300
$ws.Range("I33").Value = 300
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 300
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -71
$ws.Range("N33").ClearContents()

$ws.Range("H38").Value = 527.86664
$ws.Range("I38").Value = 221
$ws.Range("J38").Value = 878.5714
$ws.Range("K38").Value = 663
$ws.Range("L38").Value = 2635.7142
$ws.Range("M38").Value = -291
$ws.Range("N38").Value = -3379.7142

$ws.Range("H43").Value = 1854.2858
$ws.Range("I43").Value = 1745
$ws.Range("J43").Value = 1898
$ws.Range("K43").Value = 1745
$ws.Range("L43").Value = 1898
$ws.Range("M43").Value = -1676
$ws.Range("N43").Value = -2036

$ws.Range("H58").Value = 3245.5
$ws.Range("I58").Value = 242.5
$ws.Range("J58").Value = 7750
$ws.Range("K58").Value = 727.5
$ws.Range("L58").Value = 23250
$ws.Range("M58").Value = -577.5
$ws.Range("N58").Value = -23550

$ws.Range("H69").Value = 1532
$ws.Range("I69").Value = 1900
$ws.Range("J69").Value = 1500
$ws.Range("K69").Value = 5700
$ws.Range("L69").Value = 4500
$ws.Range("M69").Value = -4826
$ws.Range("N69").Value = -6248

$ws.Range("H72").Value = 1532
$ws.Range("I72").Value = 1900
$ws.Range("J72").Value = 1500
$ws.Range("K72").Value = 17100
$ws.Range("L72").Value = 13500
$ws.Range("M72").Value = -12732
$ws.Range("N72").Value = -22236

$ws.Range("H107").Value = 807
$ws.Range("I107").Value = 805
$ws.Range("J107").Value = 811.5
$ws.Range("K107").Value = 805
$ws.Range("L107").Value = 811.5
$ws.Range("M107").Value = 1115
$ws.Range("N107").Value = -4651.5

$ws.Range("H129").Value = 741.14813
$ws.Range("I129").Value = 403.7143
$ws.Range("J129").Value = 859.25
$ws.Range("K129").Value = 1211.1429
$ws.Range("L129").Value = 2577.75
$ws.Range("M129").Value = 3788.8571
$ws.Range("N129").Value = -12577.75

$ws.Range("H132").Value = 3023.0344
$ws.Range("I132").Value = 3225.3333
$ws.Range("J132").Value = 2052
$ws.Range("K132").Value = 9675.999899999999
$ws.Range("L132").Value = 6156
$ws.Range("M132").Value = -7145.999899999999
$ws.Range("N132").Value = -11216

$ws.Range("H141").Value = 3841.3635
$ws.Range("I141").Value = 2830
$ws.Range("J141").Value = 4684.1665
$ws.Range("K141").Value = 8490
$ws.Range("L141").Value = 14052.4995
$ws.Range("M141").Value = -3310
$ws.Range("N141").Value = -24412.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 72.5
$ws.Range("I4").Value = 72.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 72.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 43.5

$ws.Range("H32").Value = 4056.8845
$ws.Range("I32").Value = 3567.2727
$ws.Range("J32").Value = 6749.75
$ws.Range("K32").Value = 3567.2727
$ws.Range("L32").Value = 6749.75
$ws.Range("M32").Value = -3280.2727
$ws.Range("N32").Value = -7323.75

$ws.Range("H61").Value = 2512.8484
$ws.Range("I61").Value = 1371.7778
$ws.Range("J61").Value = 3882.1333
$ws.Range("K61").Value = 1371.7778
$ws.Range("L61").Value = 3882.1333
$ws.Range("M61").Value = -1159.7778
$ws.Range("N61").Value = -4306.1333

$ws.Range("H122").Value = 2295.1538
$ws.Range("I122").Value = 2447
$ws.Range("J122").Value = 1953.5
$ws.Range("K122").Value = 7341
$ws.Range("L122").Value = 5860.5
$ws.Range("M122").Value = -4891
$ws.Range("N122").Value = -10760.5

$ws.Range("H136").Value = 2512.8484
$ws.Range("I136").Value = 1371.7778
$ws.Range("J136").Value = 3882.1333
$ws.Range("K136").Value = 4115.3334
$ws.Range("L136").Value = 11646.3999
$ws.Range("M136").Value = -1565.3334
$ws.Range("N136").Value = -16746.3999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1649.5714
$ws.Range("I20").Value = 1771.3125
$ws.Range("J20").Value = 1260
$ws.Range("K20").Value = 1771.3125
$ws.Range("L20").Value = 1260
$ws.Range("M20").Value = -1524.3125
$ws.Range("N20").Value = -1754

$ws.Range("H22").Value = 284.25
$ws.Range("I22").Value = 284.25
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 284.25
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -111.25
$ws.Range("N22").ClearContents()

$ws.Range("H86").Value = 1829.5714
$ws.Range("I86").Value = 1680.6471
$ws.Range("J86").Value = 2462.5
$ws.Range("K86").Value = 1680.6471
$ws.Range("L86").Value = 2462.5
$ws.Range("M86").Value = -557.6470999999999
$ws.Range("N86").Value = -4708.5

$ws.Range("H89").Value = 1829.5714
$ws.Range("I89").Value = 1680.6471
$ws.Range("J89").Value = 2462.5
$ws.Range("K89").Value = 8403.235499999999
$ws.Range("L89").Value = 12312.5
$ws.Range("M89").Value = -2787.235499999999
$ws.Range("N89").Value = -23544.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3809.32
$ws.Range("I31").Value = 3931.889
$ws.Range("J31").Value = 3740.375
$ws.Range("K31").Value = 3931.889
$ws.Range("L31").Value = 3740.375
$ws.Range("M31").Value = -3636.889
$ws.Range("N31").Value = -4330.375

$ws.Range("H34").Value = 3809.32
$ws.Range("I34").Value = 3931.889
$ws.Range("J34").Value = 3740.375
$ws.Range("K34").Value = 3931.889
$ws.Range("L34").Value = 3740.375
$ws.Range("M34").Value = -3729.889
$ws.Range("N34").Value = -4144.375

$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 5000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 5000
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -4251
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 5000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 15000
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -11256
$ws.Range("N72").ClearContents()

$ws.Range("H134").Value = 1194.8695
$ws.Range("I134").Value = 897.5
$ws.Range("J134").Value = 1657.4445
$ws.Range("K134").Value = 2692.5
$ws.Range("L134").Value = 4972.333500000001
$ws.Range("M134").Value = -157.5
$ws.Range("N134").Value = -10042.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 300
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -524

$ws.Range("H5").Value = 1529.5834
$ws.Range("I5").Value = 1395.909
$ws.Range("J5").Value = 3000
$ws.Range("K5").Value = 4187.727000000001
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = -4075.727000000001
$ws.Range("N5").Value = -9224

$ws.Range("H62").Value = 6014.25
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 6014.25
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 18042.75
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -19414.75

$ws.Range("H64").Value = 3503.5
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3503.5
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 10510.5
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -11050.5

$ws.Range("H65").Value = 6014.25
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 6014.25
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 54128.25
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -60992.25

$ws.Range("H67").Value = 3503.5
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3503.5
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 10510.5
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -12382.5

$ws.Range("H74").Value = 8724.75
$ws.Range("I74").Value = 5200
$ws.Range("J74").Value = 9899.666999999999
$ws.Range("K74").Value = 15600
$ws.Range("L74").Value = 29699.001
$ws.Range("M74").Value = -14539
$ws.Range("N74").Value = -31821.001

$ws.Range("H77").Value = 8724.75
$ws.Range("I77").Value = 5200
$ws.Range("J77").Value = 9899.666999999999
$ws.Range("K77").Value = 46800
$ws.Range("L77").Value = 89097.003
$ws.Range("M77").Value = -41496
$ws.Range("N77").Value = -99705.003

$ws.Range("H80").Value = 2899.75
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2899.75
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 8699.25
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -10571.25

$ws.Range("H83").Value = 2899.75
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2899.75
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 26097.75
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -35457.75

$ws.Range("H131").Value = 792.95
$ws.Range("I131").Value = 515
$ws.Range("J131").Value = 798.62244
$ws.Range("K131").Value = 1545
$ws.Range("L131").Value = 2395.86732
$ws.Range("M131").Value = 3495
$ws.Range("N131").Value = -12475.86732

$ws.Range("H135").Value = 1529.5834
$ws.Range("I135").Value = 1395.909
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 12563.181
$ws.Range("L135").Value = 27000
$ws.Range("M135").Value = -10028.181
$ws.Range("N135").Value = -32070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6263470
$ws.Range("I70").Value = 49900
$ws.Range("J70").Value = 7816862.5
$ws.Range("K70").Value = 49900
$ws.Range("L70").Value = 7816862.5
$ws.Range("M70").Value = -49630
$ws.Range("N70").Value = -7817402.5

$ws.Range("H73").Value = 6263470
$ws.Range("I73").Value = 49900
$ws.Range("J73").Value = 7816862.5
$ws.Range("K73").Value = 49900
$ws.Range("L73").Value = 7816862.5
$ws.Range("M73").Value = -48964
$ws.Range("N73").Value = -7818734.5

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H102").Value = 1571.8422
$ws.Range("I102").Value = 1626.5625
$ws.Range("J102").Value = 1280
$ws.Range("K102").Value = 1626.5625
$ws.Range("L102").Value = 1280
$ws.Range("M102").Value = -4.5625
$ws.Range("N102").Value = -4524

$ws.Range("H122").Value = 2769.3333
$ws.Range("I122").Value = 2250
$ws.Range("J122").Value = 3184.8
$ws.Range("K122").Value = 6750
$ws.Range("L122").Value = 9554.400000000001
$ws.Range("M122").Value = -4300
$ws.Range("N122").Value = -14454.4

$ws.Range("H132").Value = 31664.834
$ws.Range("I132").Value = 4179.1816
$ws.Range("J132").Value = 74856.57000000001
$ws.Range("K132").Value = 12537.5448
$ws.Range("L132").Value = 224569.71
$ws.Range("M132").Value = -10007.5448
$ws.Range("N132").Value = -229629.71

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3224.1428
$ws.Range("I7").Value = 3429.1428
$ws.Range("J7").Value = 2814.1428
$ws.Range("K7").Value = 3429.1428
$ws.Range("L7").Value = 2814.1428
$ws.Range("M7").Value = -3317.1428
$ws.Range("N7").Value = -3038.1428

$ws.Range("H93").Value = 4375
$ws.Range("I93").Value = 4375
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 4375
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -3127

$ws.Range("H126").Value = 3224.1428
$ws.Range("I126").Value = 3429.1428
$ws.Range("J126").Value = 2814.1428
$ws.Range("K126").Value = 10287.4284
$ws.Range("L126").Value = 8442.428400000001
$ws.Range("M126").Value = -7817.428400000001
$ws.Range("N126").Value = -13382.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 17858324
$ws.Range("I136").Value = 25001014
$ws.Range("J136").Value = 1600.1875
$ws.Range("K136").Value = 75003042
$ws.Range("L136").Value = 4800.5625
$ws.Range("M136").Value = -75000492
$ws.Range("N136").Value = -9900.5625
